$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 5.763841670213395
$ws.Range("A3").Value = 9.126103885849659
$ws.Range("A4").Value = 10.36343549267141
$ws.Range("A5").Value = 11.40455399623508
$ws.Range("A6").Value = 6.035988099228689
$ws.Range("A7").Value = 11.68693180418234
$ws.Range("A8").Value = 8.589748610734119
$ws.Range("A9").Value = 7.436663384409371
$ws.Range("A10").Value = 7.815448069575837
$ws.Range("A11").Value = 7.945335257886313
$ws.Range("A12").Value = 2.451116197982941
$ws.Range("A13").Value = 6.178423552399238
$ws.Range("A14").Value = 4.518903128704466
$ws.Range("A15").Value = 2.100334423147075
$ws.Range("A16").Value = 2.343551991077874

# Rows A17:A86 all share the same updated value
$ws.Range("A17:A86").Value = 5.064110380611197
